$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of projected_week (column B) -> new projected_kcp value (column D)
$weekValues = @{
    1 = 0.1093562
    2 = 0.1099383
    3 = 0.1105753
    4 = 0.1111836
    5 = 0.1116654
    6 = 0.1119481
    7 = 0.1120475
    8 = 0.1121255
    9 = 0.1124935
    10 = 0.1135478
    11 = 0.1156918
    12 = 0.1193118
    13 = 0.1248228
    14 = 0.1327456
    15 = 0.1437692
    16 = 0.1587521
    17 = 0.1786048
    18 = 0.2039923
    19 = 0.2348632
    20 = 0.2699906
    21 = 0.3068962
    22 = 0.3424243
    23 = 0.3737478
    24 = 0.3991903
    25 = 0.4184241
    26 = 0.432119
    27 = 0.4414146
    28 = 0.4475087
    29 = 0.4514434
    30 = 0.4540282
    31 = 0.4558202
    32 = 0.4571062
    33 = 0.4578825
    34 = 0.4578512
    35 = 0.4564597
    36 = 0.452982
    37 = 0.4465912
    38 = 0.4363407
    39 = 0.4210124
    40 = 0.3989418
    41 = 0.368237
    42 = 0.3280665
    43 = 0.280919
    44 = 0.2333281
    45 = 0.1925555
    46 = 0.1622489
    47 = 0.1418359
    48 = 0.1287973
    49 = 0.1206138
    50 = 0.1154491
    51 = 0.1121294
    52 = 0.1098268
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $week = [int]$ws.Cells.Item($r, 2).Value()
    $newVal = $weekValues[$week]
    $ws.Cells.Item($r, 4).Value = $newVal
}
